$d = $word.ActiveDocument

$find = $d.Content.Find
$find.Execute("sense", $true, $false, $false, $false, $false, $true, 1, $false, "sense,", 2)
